# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Betarraga" at row 213, pushing the
# existing rows 213-264 down to 214-265 (this also extends the sheet
# dimension from A1:R264 to A1:R265 automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(213).Insert()

$ws.Cells.Item(213, 1).Value = 7
$ws.Cells.Item(213, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(213, 3).Value = "Ñuble"
$ws.Cells.Item(213, 4).Value = 44476
$ws.Cells.Item(213, 5).Value = 16
$ws.Cells.Item(213, 6).Value = 100114014
$ws.Cells.Item(213, 7).Value = "Betarraga"
$ws.Cells.Item(213, 8).Value = "Sin especificar"
$ws.Cells.Item(213, 9).Value = "Primera"
$ws.Cells.Item(213, 10).Value = 160
$ws.Cells.Item(213, 11).Value = 700
$ws.Cells.Item(213, 12).Value = 750
$ws.Cells.Item(213, 13).Value = 725
$ws.Cells.Item(213, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(213, 15).Value = "Región del Maule"
$ws.Cells.Item(213, 16).Value = 145
$ws.Cells.Item(213, 17).Value = 5
$ws.Cells.Item(213, 18).Value = "Hortaliza"
